$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the very top. This shifts every existing row
# (including the old text header in row 1) down by one.
$ws.Rows.Item(1).Insert()

# The old header row (text labels) is now row 2, and it kept its bold /
# bordered / centered header style (style index "1"). Grab that formatting
# and stamp it onto the new row 1 before we overwrite row 2's style.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)   # xlPasteFormats

# New row 1 becomes a simple numeric column-index row: 0, 1, 2, ... 11
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# Row 2 (the old header text) loses the special header styling - it goes
# back to the default/normal style, matching the rest of the data rows.
$ws.Rows.Item(2).Style = "Normal"

# Row 2's I/K/L cells are blanked out (the old K1/L1 "thread_size" /
# "material_surface" helper labels are dropped, and I2 stays blank like
# the original I1 was).
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
